$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the cells that would otherwise be
# auto-converted to a date / number by Excel's input parsing, then
# write the values as plain strings.
$ws.Range("A1:K1").NumberFormat = "@"

$ws.Range("A1").Value = "DePaola"
$ws.Range("B1").Value = "Andrew"
$ws.Range("C1").Value = "WR"
$ws.Range("D1").Value = "2018-09-10"
$ws.Range("E1").Value = "1"
$ws.Range("F1").Value = "31.044"
$ws.Range("G1").Value = "OAK"
$ws.Range("H1").Value = ""
$ws.Range("I1").Value = "LAR"
$ws.Range("J1").Value = "L 13-33"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = 0

# Drop the temporary text formatting again so the saved styles.xml
# matches the original (unstyled) cell formatting.
$ws.Range("A1:L1").ClearFormats()
